$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. The "Climate modes ... beyond" sentence and the following "."
#    currently live in two separate runs. Run a Find/Replace whose
#    match spans both runs (replacing the tail of the sentence plus the
#    period with the same literal text) so Word collapses them into a
#    single run containing the whole sentence, ending in the period.
# ---------------------------------------------------------------------
$sentenceTail = "time scales varying from days to seasons and beyond."
$d.Content.Find.Execute(
    $sentenceTail, $true, $false, $false, $false, $false,
    $true, 1, $false, $sentenceTail, 2) | Out-Null

$fullSentence = "Climate modes are recurrent patterns, usually of pressure or sea-surface temperature (SST), typically characterized by negative and positive phases, which each have distinctive effects on the distribution of rainfall, temperature and other meteorological elements on time scales varying from days to seasons and beyond."

# Locate that paragraph (its full text now matches exactly; Range.Text
# includes the trailing paragraph-mark character, so trim before
# comparing).
$climateParaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.TrimEnd() -eq $fullSentence) {
        $climateParaIndex = $i
        break
    }
}

if ($climateParaIndex -eq -1) {
    throw "Could not locate the 'Climate modes...' paragraph"
}

$climatePara = $d.Paragraphs($climateParaIndex)

# ---------------------------------------------------------------------
# 2. Insert a new Heading1 paragraph "precip_quantiles_9month" right
#    after the "Climate modes..." paragraph.
# ---------------------------------------------------------------------
$climatePara.Range.InsertParagraphAfter() | Out-Null
$headingPara = $d.Paragraphs($climateParaIndex + 1)
$headingPara.Range.Text = "precip_quantiles_9month"
$headingPara.Style = "Heading1"

# ---------------------------------------------------------------------
# 3. Insert a new body paragraph after the heading, describing the
#    precip_quantiles_9month product.
# ---------------------------------------------------------------------
$headingPara.Range.InsertParagraphAfter() | Out-Null
$bodyPara = $d.Paragraphs($climateParaIndex + 2)
$bodyPara.Style = "Normal"
$bodyPara.Range.Text = "Precipitation quantiles are based on the nine months aggregated GPCC Monitoring Product and First Guess Monthly product. The baseline period is 1991-2020, using Full Data Monthly in its latest version. Quality controlled rain gauge (in situ) data are used and the quality control protocol depends on the data set."
